# Auto-generated edit script applying diff changes to Durandal_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2747.2222
$ws.Cells.Item(86, 9).Value = 1579.5714
$ws.Cells.Item(86, 10).Value = 6834
$ws.Cells.Item(86, 11).Value = 1579.5714
$ws.Cells.Item(86, 12).Value = 6834
$ws.Cells.Item(86, 13).Value = -456.5714
$ws.Cells.Item(86, 14).Value = -9080

$ws.Cells.Item(88, 8).Value = 6814030.5
$ws.Cells.Item(88, 10).Value = 9957475
$ws.Cells.Item(88, 12).Value = 9957475
$ws.Cells.Item(88, 14).Value = -9958287

$ws.Cells.Item(89, 8).Value = 2747.2222
$ws.Cells.Item(89, 9).Value = 1579.5714
$ws.Cells.Item(89, 10).Value = 6834
$ws.Cells.Item(89, 11).Value = 7897.857
$ws.Cells.Item(89, 12).Value = 34170
$ws.Cells.Item(89, 13).Value = -2281.857
$ws.Cells.Item(89, 14).Value = -45402

$ws.Cells.Item(91, 8).Value = 6814030.5
$ws.Cells.Item(91, 10).Value = 9957475
$ws.Cells.Item(91, 12).Value = 9957475
$ws.Cells.Item(91, 14).Value = -9960283

$ws.Cells.Item(129, 8).Value = 1003.1724
$ws.Cells.Item(129, 9).Value = 339.4
$ws.Cells.Item(129, 10).Value = 1141.4584
$ws.Cells.Item(129, 11).Value = 1018.2
$ws.Cells.Item(129, 12).Value = 3424.3752
$ws.Cells.Item(129, 13).Value = 3981.8
$ws.Cells.Item(129, 14).Value = -13424.3752

$ws.Cells.Item(138, 8).Value = 3276.59
$ws.Cells.Item(138, 9).Value = 2263.7856
$ws.Cells.Item(138, 10).Value = 3578.2766
$ws.Cells.Item(138, 11).Value = 6791.3568
$ws.Cells.Item(138, 12).Value = 10734.8298
$ws.Cells.Item(138, 13).Value = -1651.3568
$ws.Cells.Item(138, 14).Value = -21014.8298

$ws.Cells.Item(139, 8).Value = 39820
$ws.Cells.Item(139, 10).Value = 39820
$ws.Cells.Item(139, 12).Value = 39820
$ws.Cells.Item(139, 14).Value = -50100

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 383198.03
$ws.Cells.Item(32, 9).Value = 4092
$ws.Cells.Item(32, 11).Value = 4092
$ws.Cells.Item(32, 13).Value = -3805

$ws.Cells.Item(110, 8).Value = 930.0833
$ws.Cells.Item(110, 9).Value = 832.8182
$ws.Cells.Item(110, 10).Value = 2000
$ws.Cells.Item(110, 11).Value = 832.8182
$ws.Cells.Item(110, 12).Value = 2000
$ws.Cells.Item(110, 13).Value = 1212.1818
$ws.Cells.Item(110, 14).Value = -6090

$ws.Cells.Item(132, 8).Value = 18888742
$ws.Cells.Item(132, 9).Value = 28573036
$ws.Cells.Item(132, 10).Value = 58170.445
$ws.Cells.Item(132, 11).Value = 85719108
$ws.Cells.Item(132, 12).Value = 174511.335
$ws.Cells.Item(132, 13).Value = -85716578
$ws.Cells.Item(132, 14).Value = -179571.335

$ws.Cells.Item(133, 8).Value = 31800
$ws.Cells.Item(133, 10).Value = 31800
$ws.Cells.Item(133, 12).Value = 31800
$ws.Cells.Item(133, 14).Value = -36860

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 8942.5
$ws.Cells.Item(20, 9).Value = 8199.546
$ws.Cells.Item(20, 10).Value = 11666.667
$ws.Cells.Item(20, 11).Value = 8199.546
$ws.Cells.Item(20, 12).Value = 11666.667
$ws.Cells.Item(20, 13).Value = -7952.546
$ws.Cells.Item(20, 14).Value = -12160.667

$ws.Cells.Item(86, 8).Value = 35752090
$ws.Cells.Item(86, 9).Value = 71430890
$ws.Cells.Item(86, 10).Value = 73288.14
$ws.Cells.Item(86, 11).Value = 71430890
$ws.Cells.Item(86, 12).Value = 73288.14
$ws.Cells.Item(86, 13).Value = -71429767
$ws.Cells.Item(86, 14).Value = -75534.14

$ws.Cells.Item(89, 8).Value = 35752090
$ws.Cells.Item(89, 9).Value = 71430890
$ws.Cells.Item(89, 10).Value = 73288.14
$ws.Cells.Item(89, 11).Value = 357154450
$ws.Cells.Item(89, 12).Value = 366440.7
$ws.Cells.Item(89, 13).Value = -357148834
$ws.Cells.Item(89, 14).Value = -377672.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(33, 8).Value = 4700
$ws.Cells.Item(33, 10).Value = 4700
$ws.Cells.Item(33, 12).Value = 4700
$ws.Cells.Item(33, 14).Value = -5458

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 16818
$ws.Cells.Item(3, 9).Value = 16818
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 50454
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).Value = -50342

$ws.Cells.Item(75, 8).Value = 866.6667
$ws.Cells.Item(75, 9).Value = 600
$ws.Cells.Item(75, 10).Value = 1000
$ws.Cells.Item(75, 11).Value = 1800
$ws.Cells.Item(75, 12).Value = 3000
$ws.Cells.Item(75, 13).Value = -802
$ws.Cells.Item(75, 14).Value = -4996

$ws.Cells.Item(78, 8).Value = 866.6667
$ws.Cells.Item(78, 9).Value = 600
$ws.Cells.Item(78, 10).Value = 1000
$ws.Cells.Item(78, 11).Value = 5400
$ws.Cells.Item(78, 12).Value = 9000
$ws.Cells.Item(78, 13).Value = -408
$ws.Cells.Item(78, 14).Value = -18984

$ws.Cells.Item(92, 8).Value = 514.8570999999999
$ws.Cells.Item(92, 9).Value = 650
$ws.Cells.Item(92, 10).Value = 460.8
$ws.Cells.Item(92, 11).Value = 1950
$ws.Cells.Item(92, 12).Value = 1382.4
$ws.Cells.Item(92, 13).Value = -702
$ws.Cells.Item(92, 14).Value = -3878.4

$ws.Cells.Item(113, 8).Value = 1009.58905
$ws.Cells.Item(113, 9).Value = 826.5
$ws.Cells.Item(113, 10).Value = 1014.74646
$ws.Cells.Item(113, 11).Value = 2479.5
$ws.Cells.Item(113, 12).Value = 3044.23938
$ws.Cells.Item(113, 13).Value = -309.5
$ws.Cells.Item(113, 14).Value = -7384.23938

$ws.Cells.Item(121, 8).Value = 31775.363
$ws.Cells.Item(121, 9).Value = 6216
$ws.Cells.Item(121, 10).Value = 36339.535
$ws.Cells.Item(121, 11).Value = 18648
$ws.Cells.Item(121, 12).Value = 109018.605
$ws.Cells.Item(121, 13).Value = -17338
$ws.Cells.Item(121, 14).Value = -111638.605

$ws.Cells.Item(122, 8).Value = 657.9259
$ws.Cells.Item(122, 9).Value = 298.4762
$ws.Cells.Item(122, 10).Value = 1916
$ws.Cells.Item(122, 11).Value = 2686.2858
$ws.Cells.Item(122, 12).Value = 17244
$ws.Cells.Item(122, 13).Value = -236.2858000000001
$ws.Cells.Item(122, 14).Value = -22144

$ws.Cells.Item(134, 8).Value = 5439.9707
$ws.Cells.Item(134, 9).Value = 5866.263
$ws.Cells.Item(134, 11).Value = 17598.789
$ws.Cells.Item(134, 13).Value = -12528.789

$ws.Cells.Item(138, 8).Value = 1387.2727
$ws.Cells.Item(138, 9).Value = 607.5
$ws.Cells.Item(138, 10).Value = 3466.6667
$ws.Cells.Item(138, 11).Value = 1822.5
$ws.Cells.Item(138, 12).Value = 10400.0001
$ws.Cells.Item(138, 13).Value = 3317.5
$ws.Cells.Item(138, 14).Value = -20680.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 1239.6666
$ws.Cells.Item(36, 9).Value = 1500
$ws.Cells.Item(36, 10).Value = 1187.6
$ws.Cells.Item(36, 11).Value = 1500
$ws.Cells.Item(36, 12).Value = 1187.6
$ws.Cells.Item(36, 13).Value = -1015
$ws.Cells.Item(36, 14).Value = -2157.6

$ws.Cells.Item(43, 8).Value = 22539.666
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 22539.666
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).ClearContents()
$ws.Cells.Item(43, 13).Value = 22539.666
$ws.Cells.Item(43, 14).Value = -22841.666

$ws.Cells.Item(46, 8).Value = 12660
$ws.Cells.Item(46, 9).Value = 3300
$ws.Cells.Item(46, 10).Value = 22020
$ws.Cells.Item(46, 11).Value = 3300
$ws.Cells.Item(46, 12).Value = 22020
$ws.Cells.Item(46, 13).Value = -3144
$ws.Cells.Item(46, 14).Value = -22332

$ws.Cells.Item(102, 8).Value = 4999
$ws.Cells.Item(102, 10).Value = 4999
$ws.Cells.Item(102, 12).Value = 4999
$ws.Cells.Item(102, 14).Value = -8243

$ws.Cells.Item(126, 8).Value = 7581636.5
$ws.Cells.Item(126, 9).Value = 8819.923000000001
$ws.Cells.Item(126, 10).Value = 18520150
$ws.Cells.Item(126, 11).Value = 26459.769
$ws.Cells.Item(126, 12).Value = 55560450
$ws.Cells.Item(126, 13).Value = -23989.769
$ws.Cells.Item(126, 14).Value = -55565390

$ws.Cells.Item(132, 8).Value = 269571.12
$ws.Cells.Item(132, 9).Value = 34179.2
$ws.Cells.Item(132, 10).Value = 911549.0600000001
$ws.Cells.Item(132, 11).Value = 102537.6
$ws.Cells.Item(132, 12).Value = 2734647.18
$ws.Cells.Item(132, 13).Value = -100007.6
$ws.Cells.Item(132, 14).Value = -2739707.18

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3021.9412
$ws.Cells.Item(7, 9).Value = 1999.75
$ws.Cells.Item(7, 10).Value = 3336.4614
$ws.Cells.Item(7, 11).Value = 1999.75
$ws.Cells.Item(7, 12).Value = 3336.4614
$ws.Cells.Item(7, 13).Value = -1887.75
$ws.Cells.Item(7, 14).Value = -3560.4614

$ws.Cells.Item(40, 8).Value = 2574.7334
$ws.Cells.Item(40, 9).Value = 2414.2856
$ws.Cells.Item(40, 10).Value = 2715.125
$ws.Cells.Item(40, 11).Value = 2414.2856
$ws.Cells.Item(40, 12).Value = 2715.125
$ws.Cells.Item(40, 13).Value = -2278.2856
$ws.Cells.Item(40, 14).Value = -2987.125

$ws.Cells.Item(68, 8).Value = 2038.6666
$ws.Cells.Item(68, 9).Value = 1941.037
$ws.Cells.Item(68, 10).Value = 2258.3333
$ws.Cells.Item(68, 11).Value = 1941.037
$ws.Cells.Item(68, 12).Value = 2258.3333
$ws.Cells.Item(68, 13).Value = -1192.037
$ws.Cells.Item(68, 14).Value = -3756.3333

$ws.Cells.Item(71, 8).Value = 2038.6666
$ws.Cells.Item(71, 9).Value = 1941.037
$ws.Cells.Item(71, 10).Value = 2258.3333
$ws.Cells.Item(71, 11).Value = 9705.184999999999
$ws.Cells.Item(71, 12).Value = 11291.6665
$ws.Cells.Item(71, 13).Value = -5961.184999999999
$ws.Cells.Item(71, 14).Value = -18779.6665

$ws.Cells.Item(82, 8).Value = 2748.8
$ws.Cells.Item(82, 9).Value = 3872
$ws.Cells.Item(82, 10).Value = 2000
$ws.Cells.Item(82, 11).Value = 3872
$ws.Cells.Item(82, 12).Value = 2000
$ws.Cells.Item(82, 13).Value = -3511
$ws.Cells.Item(82, 14).Value = -2722

$ws.Cells.Item(85, 8).Value = 2748.8
$ws.Cells.Item(85, 9).Value = 3872
$ws.Cells.Item(85, 10).Value = 2000
$ws.Cells.Item(85, 11).Value = 3872
$ws.Cells.Item(85, 12).Value = 2000
$ws.Cells.Item(85, 13).Value = -2624
$ws.Cells.Item(85, 14).Value = -4496

$ws.Cells.Item(126, 8).Value = 3021.9412
$ws.Cells.Item(126, 9).Value = 1999.75
$ws.Cells.Item(126, 10).Value = 3336.4614
$ws.Cells.Item(126, 11).Value = 5999.25
$ws.Cells.Item(126, 12).Value = 10009.3842
$ws.Cells.Item(126, 13).Value = -3529.25
$ws.Cells.Item(126, 14).Value = -14949.3842

$ws.Cells.Item(132, 8).Value = 24325
$ws.Cells.Item(132, 9).Value = 40020.348
$ws.Cells.Item(132, 10).Value = 1653.9445
$ws.Cells.Item(132, 11).Value = 120061.044
$ws.Cells.Item(132, 12).Value = 4961.833500000001
$ws.Cells.Item(132, 13).Value = -117531.044
$ws.Cells.Item(132, 14).Value = -10021.8335

$ws.Cells.Item(133, 8).Value = 111181.82
$ws.Cells.Item(133, 10).Value = 111181.82
$ws.Cells.Item(133, 12).Value = 111181.82
$ws.Cells.Item(133, 14).Value = -116241.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 65497610
$ws.Cells.Item(132, 9).Value = 86924424
$ws.Cells.Item(132, 10).Value = 3597929
$ws.Cells.Item(132, 11).Value = 260773272
$ws.Cells.Item(132, 12).Value = 10793787
$ws.Cells.Item(132, 13).Value = -260770742
$ws.Cells.Item(132, 14).Value = -10798847
